$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) with new short English column names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize Spanish connector words (de, del, la, las, los, el, y) to title case ---
# --- in municipality / state name cells (and fix "MonteMorelos" typo)            ---
$renames = @(
    @('B7', 'Pabellón De Arteaga'),
    @('B8', 'Rincón De Romos'),
    @('B9', 'San Francisco De Los Romo'),
    @('B41', 'Mazapa De Madero'),
    @('B42', 'Montecristo De Guerrero'),
    @('B69', 'Guadalupe Y Calvo'),
    @('B70', 'Hidalgo Del Parral'),
    @('B76', 'San Francisco Del Oro'),
    @('B78', 'Valle De Zaragoza'),
    @('B86', 'San Juan De Sabinas'),
    @('B93', 'Villa De Álvarez'),
    @('A95', 'Ciudad De México'),
    @('B99', 'Cuajimalpa De Morelos'),
    @('B121', 'Nombre De Dios'),
    @('B128', 'San Juan De Guadalupe'),
    @('B129', 'San Juan Del Río'),
    @('A136', 'Estado De México'),
    @('B136', 'Acambay De Ruíz Castañeda'),
    @('B144', 'Chapa De Mota'),
    @('B148', 'Ecatepec De Morelos'),
    @('B151', 'Ixtapan De La Sal'),
    @('B157', 'Naucalpan De Juárez'),
    @('B163', 'San Antonio La Isla'),
    @('B164', 'San Felipe Del Progreso'),
    @('B168', 'Tenango Del Valle'),
    @('B169', 'Tlalnepantla De Baz'),
    @('B173', 'Villa De Allende'),
    @('B182', 'Apaseo El Alto'),
    @('B183', 'Apaseo El Grande'),
    @('B189', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B193', 'Jaral Del Progreso'),
    @('B198', 'Purísima Del Rincón'),
    @('B202', 'San Diego De La Unión'),
    @('B204', 'San Francisco Del Rincón'),
    @('B206', 'San Luis De La Paz'),
    @('B207', 'Santa Cruz De Juventino Rosas'),
    @('B208', 'Silao De La Victoria'),
    @('B212', 'Valle De Santiago'),
    @('B218', 'Acapulco De Juárez'),
    @('B219', 'Ajuchitlán Del Progreso'),
    @('B222', 'Atlamajalcingo Del Monte'),
    @('B224', 'Atoyac De Álvarez'),
    @('B225', 'Ayutla De Los Libres'),
    @('B227', 'Chilapa De Álvarez'),
    @('B228', 'Chilpancingo De Los Bravo'),
    @('B229', 'Coahuayutla De José María Izazaga'),
    @('B232', 'Coyuca De Benítez'),
    @('B233', 'Coyuca De Catalán'),
    @('B234', 'Cutzamala De Pinzón'),
    @('B238', 'Huitzuco De Los Figueroa'),
    @('B239', 'Iguala De La Independencia'),
    @('B241', 'Ixcateopan De Cuauhtémoc'),
    @('B243', 'La Unión De Isidoro Montes De Oca'),
    @('B245', 'Mártir De Cuilapan'),
    @('B254', 'Taxco De Alarcón'),
    @('B256', 'Técpan De Galeana'),
    @('B258', 'Tepecoacuilco De Trujano'),
    @('B262', 'Tlapa De Comonfort'),
    @('B274', 'Cuautepec De Hinojosa'),
    @('B279', 'Jacala De Ledezma'),
    @('B283', 'Mineral Del Chico'),
    @('B284', 'Molango De Escamilla'),
    @('B286', 'Pachuca De Soto'),
    @('B288', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B292', 'Tenango De Doria'),
    @('B293', 'Tepeji Del Río De Ocampo'),
    @('B297', 'Tula De Allende'),
    @('B298', 'Tulancingo De Bravo'),
    @('B299', 'Zacualtipán De Ángeles'),
    @('B305', 'Atemajac De Brizuela'),
    @('B306', 'Atotonilco El Alto'),
    @('B308', 'Autlán De Navarro'),
    @('B314', 'Encarnación De Díaz'),
    @('B318', 'Huejuquilla El Alto'),
    @('B319', 'Ixtlahuacán Del Río'),
    @('B324', 'Lagos De Moreno'),
    @('B327', 'Ojuelos De Jalisco'),
    @('B329', 'San Diego De Alejandría'),
    @('B336', 'Teocuitatlán De Corona'),
    @('B337', 'Tepatitlán De Morelos'),
    @('B339', 'Tizapán El Alto'),
    @('B344', 'Unión De San Antonio'),
    @('B347', 'Yahualica De González Gallo'),
    @('B348', 'Zacoalco De Torres'),
    @('B351', 'Zapotlán El Grande'),
    @('B364', 'Coalcomán De Vázquez Pallares'),
    @('B398', 'Tiquicheo De Nicolás Romero'),
    @('B423', 'Tlaltizapán De Zapata'),
    @('B429', 'Amatlán De Cañas'),
    @('B445', 'Montemorelos'),
    @('B448', 'San Nicolás De Los Garza'),
    @('B450', 'Acatlán De Pérez Figueroa'),
    @('B455', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B456', 'Heroica Ciudad De Huajuapan De León'),
    @('B457', 'Heroica Ciudad De Tlaxiaco'),
    @('B459', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B463', 'Miahuatlán De Porfirio Díaz'),
    @('B465', 'Oaxaca De Juárez'),
    @('B466', 'Ocotlán De Morelos'),
    @('B521', 'Villa De Tututepec De Melchor Ocampo'),
    @('B522', 'Villa Sola De Vega'),
    @('B542', 'Ixcamilpa De Guerrero'),
    @('B545', 'Izúcar De Matamoros'),
    @('B553', 'Tepanco De López'),
    @('B563', 'Xochitlán De Vicente Suárez'),
    @('B570', 'Amealco De Bonfil'),
    @('B574', 'Jalpan De Serra'),
    @('B575', 'Landa De Matamoros'),
    @('B577', 'Pinal De Amoles'),
    @('B579', 'San Juan Del Río'),
    @('B587', 'Armadillo De Los Infante'),
    @('B591', 'Ciudad Del Maíz'),
    @('B596', 'Mexquitic De Carmona'),
    @('B601', 'San Ciro De Acosta'),
    @('B605', 'Santa María Del Río'),
    @('B611', 'Tanquián De Escobedo'),
    @('B614', 'Villa De Arista'),
    @('B615', 'Villa De Arriaga'),
    @('B616', 'Villa De Guadalupe'),
    @('B617', 'Villa De Ramos'),
    @('B618', 'Villa De Reyes'),
    @('B665', 'Soto La Marina'),
    @('B671', 'Acuamanala De Miguel Hidalgo'),
    @('B674', 'Ixtacuixtla De Mariano Matamoros'),
    @('B676', 'Papalotla De Xicohténcatl'),
    @('B683', 'Amatlán De Los Reyes'),
    @('B688', 'Boca Del Río'),
    @('B696', 'Cosamaloapan De Carpio'),
    @('B702', 'Ixhuatlán De Madero'),
    @('B708', 'Martínez De La Torre'),
    @('B719', 'Poza Rica De Hidalgo'),
    @('B724', 'Soledad De Doblado'),
    @('B756', 'Nochistlán De Mejía'),
    @('B757', 'Noria De Ángeles'),
    @('B766', 'Tlaltenango De Sánchez Román'),
    @('B768', 'Villa De Cos'),
)

foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Remove trailing footnote/metadata rows (776-780) ---
$ws.Range("A776:A780").EntireRow.Delete() | Out-Null

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()